$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 118, shifting existing rows 118.. down by one.
$ws.Rows.Item(118).Insert()

# Fill in the new row's data.
$ws.Cells.Item(118, 1).Value = "Nonviolent & Violent Campaigns and Outcomes"
$ws.Cells.Item(118, 2).Value = "international relations"
$ws.Cells.Item(118, 3).Value = "https://www.du.edu/korbel/sie/research/chenow_navco_data.html"
$ws.Cells.Item(118, 4).Value = "Nonviolent campaigns, violent campaigns"
$ws.Cells.Item(118, 5).Value = "world"
$ws.Cells.Item(118, 6).Value = 1900
$ws.Cells.Item(118, 7).Value = 2011
$ws.Cells.Item(118, 8).Value = "online"
$ws.Cells.Item(118, 9).Value = "free, online registration"

# The link cell (column C) gets a hyperlink + the "Hyperlink" cell style.
$ws.Hyperlinks.Add($ws.Cells.Item(118, 3), "https://www.du.edu/korbel/sie/research/chenow_navco_data.html")
$ws.Cells.Item(118, 3).Style = "Hyperlink"

# Keep the view/selection roughly where the diff shows it ended up.
$ws.Application.ActiveWindow.ScrollRow = 99
$ws.Range("A118").Select()
